$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

# The sheet's columns D:G hold (in this order):
#   D = codeforiati:group-name   E = codeforiati:group-code
#   F = codeforiati:category-code  G = codeforiati:category-name
# They need to be re-ordered to:
#   D = codeforiati:group-code   E = codeforiati:category-name
#   F = codeforiati:group-name   G = codeforiati:category-code
# i.e. newD=oldE, newE=oldG, newF=oldD, newG=oldF.
# Use Range.Copy (not value assignment) so the shared-string cell type
# and "no explicit style" state of each cell is preserved verbatim,
# rather than Excel re-inferring a numeric type / adding a style record
# for text that looks numeric (e.g. "110").

$stageD = $ws.Range("I1:I$lastRow")
$stageE = $ws.Range("J1:J$lastRow")
$stageF = $ws.Range("K1:K$lastRow")
$stageG = $ws.Range("L1:L$lastRow")

$ws.Range("D1:D$lastRow").Copy($stageD)
$ws.Range("E1:E$lastRow").Copy($stageE)
$ws.Range("F1:F$lastRow").Copy($stageF)
$ws.Range("G1:G$lastRow").Copy($stageG)

$stageE.Copy($ws.Range("D1:D$lastRow"))
$stageG.Copy($ws.Range("E1:E$lastRow"))
$stageD.Copy($ws.Range("F1:F$lastRow"))
$stageF.Copy($ws.Range("G1:G$lastRow"))

$ws.Range("I1:L$lastRow").ClearContents()
